# Updated symbol list (Price / Volume(1h) columns) per the Wed Jan 25 14:47:26 UTC 2023
# GitHub Actions data refresh. Source cells are plain text (e.g. "300.60", "-4.37%"),
# not real numbers/percentages, so each cell is force-formatted as Text before the
# new value is written (otherwise Excel would auto-convert the literal into a number)
# and then restored to the default "Normal" style so no stray formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "D2" '300.83'
Set-TextValue "E2" '-4.42%'
Set-TextValue "D3" '35.37'
Set-TextValue "E3" '-0.78%'
Set-TextValue "D4" '5.037'
Set-TextValue "E4" '-1.22%'
Set-TextValue "D5" '0.07969'
Set-TextValue "E5" '-2.16%'
Set-TextValue "D6" '1.923'
Set-TextValue "E6" '-8.64%'
Set-TextValue "D7" '4.057'
Set-TextValue "E7" '-2.07%'
Set-TextValue "E8" '-1.82%'
Set-TextValue "D9" '0.9228'
Set-TextValue "E9" '-1.10%'
Set-TextValue "D10" '0.1406'
Set-TextValue "E10" '35.81%'
Set-TextValue "D11" '0.1892'
Set-TextValue "E11" '-1.30%'
Set-TextValue "D12" '0.09157'
Set-TextValue "E12" '0.32%'
Set-TextValue "D13" '0.03416'
Set-TextValue "E13" '-5.48%'
Set-TextValue "D14" '0.09857'
Set-TextValue "E14" '-0.28%'
Set-TextValue "D15" '0.001392'
Set-TextValue "E15" '-2.67%'
Set-TextValue "D16" '0.005776'
Set-TextValue "E16" '-1.07%'
Set-TextValue "D17" '3.517'
Set-TextValue "E17" '1.34%'
Set-TextValue "D18" '2.990'
Set-TextValue "E18" '0.45%'
Set-TextValue "D19" '0.3404'
Set-TextValue "E19" '-1.60%'
Set-TextValue "D20" '0.1296'
Set-TextValue "E20" '-1.33%'
Set-TextValue "D21" '5.048'
Set-TextValue "E21" '-1.37%'
Set-TextValue "D22" '0.2407'
Set-TextValue "E22" '8.55%'
Set-TextValue "D23" '0.04489'
Set-TextValue "E23" '-1.37%'
Set-TextValue "D24" '0.001217'
Set-TextValue "E24" '-2.16%'
Set-TextValue "D25" '0.004770'
Set-TextValue "E25" '-0.50%'
Set-TextValue "D26" '0.0001232'
Set-TextValue "E26" '-1.58%'
Set-TextValue "D27" '0.0003012'
Set-TextValue "E27" '-33.25%'
Set-TextValue "D39" '0.01894'
Set-TextValue "E39" '-3.43%'
Set-TextValue "E40" '-3.77%'
Set-TextValue "D41" '0.007355'
Set-TextValue "E41" '-3.34%'
Set-TextValue "D42" '0.009747'
Set-TextValue "E42" '25.55%'
Set-TextValue "D43" '0.1325'
Set-TextValue "E43" '-4.21%'
Set-TextValue "D44" '0.002044'
Set-TextValue "E44" '-2.54%'
Set-TextValue "D45" '0.009333'
Set-TextValue "E45" '-20.39%'
Set-TextValue "D46" '0.00006266'
Set-TextValue "E46" '-6.99%'
Set-TextValue "D47" '0.00000000753'
Set-TextValue "E47" '0.24%'
Set-TextValue "E48" '42.62%'
Set-TextValue "D49" '0.001665'
Set-TextValue "E49" '-2.31%'
Set-TextValue "D50" '0.00002108'
Set-TextValue "E50" '0.24%'
Set-TextValue "D51" '0.0002008'
Set-TextValue "E51" '0.24%'
